$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.893.77"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.648.53"
$ws.Range("E3").Value = "  +0.71%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "308.60"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.84%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3896"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -1.37%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3829"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -1.14%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "51.92"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +2.91%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.351"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  +0.43%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08423"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -1.19%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "23.84"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -1.23%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.075"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -0.94%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.946"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +3.45%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.00001316"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "1.648.10"
$ws.Range("E17").Value = "  +0.98%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "94.60"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.01%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06980"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  -2.86%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.931"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  +0.61%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "13.69"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").Value = "23.882.87"
$ws.Range("E24").Value = "  -0.77%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.457"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -0.01%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.947"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("E27").Value = "  -1.54%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "150.91"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -4.03%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "5.404"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +1.04%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "138.60"
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.859"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -2.80%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.523"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "1.828.62"
$ws.Range("E33").Value = "  +0.96%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.045"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +3.63%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.08028"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -1.91%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02956"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "10.96"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +3.83%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "6.650"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -0.84%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.2677"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -0.57%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.09094"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.32%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.7586"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -0.66%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "13.46"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -2.16%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.419"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.83%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "16.32"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +0.88%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.6986"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  -1.42%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.076"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -0.58%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.08267"
$cell.ClearFormats()
$ws.Range("E49").Value = "  -1.16%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "134.30"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -1.76%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.215"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -1.41%  "
